$wb = $excel.ActiveWorkbook

$newName = "294-MS-EPP-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-ONTIME"

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update the product name value on both sheets (hyphen added after "294")
$wsInput.Range("B1").Value = $newName
$wsOutput.Range("B1").Value = $newName

# Move the selection on the input sheet to B1
$wsInput.Range("B1").Select() | Out-Null

# Make the output sheet the active sheet/tab
$wsOutput.Activate()
$wsOutput.Range("B1").Select() | Out-Null
